$d = $word.ActiveDocument

# The paragraph "<id>p026r_1</id>" is currently split across three runs:
#   "<id>" (Courier New, color 7f6000) + "p026r_1" (plain) + "</id>" (Courier New, color 7f6000)
# Collapse it into a single run, merged text "<id>p026r_1</id>",
# keeping the formatting of the surrounding "<id>"/"</id>" runs.
$p = $d.Paragraphs(6)
$r = $p.Range
$r.Find.Execute("<id>p026r_1</id>", $false, $false, $false, $false, $false,
                $true, 1, $false, "<id>p026r_1</id>", 2)
